# Applies the diff: swap the match-detail columns (F:V) between several
# pairs of rows that share the same match date/time (columns A:E stay put),
# and append a brand-new match row (row 166) at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange($r1, $r2, $colStart, $colEnd) {
    $range1 = $ws.Range($colStart + $r1 + ":" + $colEnd + $r1)
    $range2 = $ws.Range($colStart + $r2 + ":" + $colEnd + $r2)
    $v1 = $range1.Value()
    $v2 = $range2.Value()
    $range1.Value = $v2
    $range2.Value = $v1
}

# Pairs of rows whose F:V (match) data must be swapped, keeping A:E (index,
# country, league, season, date) untouched on each row.
Swap-RowRange 85 86 "F" "V"
Swap-RowRange 93 94 "F" "V"
Swap-RowRange 95 96 "F" "V"
Swap-RowRange 114 115 "F" "V"
Swap-RowRange 147 148 "F" "V"

# Append the new match row 166, copying the style of row 165 first so that
# column A (index) and column E (date) keep the correct formatting.
$ws.Range("A165:V165").Copy($ws.Range("A166:V166"))

$ws.Range("A166").Value = 165
$ws.Range("B166").Value = "spain"
$ws.Range("C166").Value = "laliga2"
$ws.Range("D166").Value = "2023-2024"
$ws.Range("E166").Value = 45243.875
$ws.Range("F166").Value = "R. Oviedo"
$ws.Range("G166").Value = 1
$ws.Range("H166").Value = "FC Cartagena SAD"
$ws.Range("I166").Value = 1
$ws.Range("J166").Value = 1.76
$ws.Range("K166").Value = "06/11/2023 21:13"
$ws.Range("L166").Value = 1.85
$ws.Range("M166").Value = "13/11/2023 20:56"
$ws.Range("N166").Value = 3.47
$ws.Range("O166").Value = "06/11/2023 21:13"
$ws.Range("P166").Value = 3.25
$ws.Range("Q166").Value = "13/11/2023 20:57"
$ws.Range("R166").Value = 5.39
$ws.Range("S166").Value = "06/11/2023 21:13"
$ws.Range("T166").Value = 5.5
$ws.Range("U166").Value = "13/11/2023 20:57"
$ws.Range("V166").Value = "https://www.betexplorer.com/football/spain/laliga2/r-oviedo-fc-cartagena-sad/dWwodK0g/"
